$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, pushing existing rows 26:168 down to 27:169
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new data point.
# Metadata columns (A,B,C,E,F,G,H,I,N,O,Q,R) are constant across the whole sheet,
# so reuse the same values found elsewhere (e.g. row 27, which now holds the
# record that used to be in row 26 before the insert).
$ws.Range("A26").Value = 8
$ws.Range("B26").Value = "Terminal La Palmera de La Serena"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44819
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = 100112040
$ws.Range("G26").Value = "Cilantro"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 2400
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 2500
$ws.Range("M26").Value = 2250
$ws.Range("N26").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O26").Value = "Provincia del Elquí"
$ws.Range("P26").Value = 1500
$ws.Range("Q26").Value = 1.5
$ws.Range("R26").Value = "Hortaliza"
